$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new data row (17) for period "2508", right below the
#        existing "2507" row (16) and above the signature block, which
#        shifts down from rows 21-22 to rows 22-23. ---
$ws.Rows("17").Insert()

# Copy the formatting (borders/fill/font/alignment) of row 16 onto the new
# row 17 so it matches the rest of the data table.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's values (same worker, new "2508" period).
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1007958543"
$ws.Range("D17").Value = "JESUS DAVID MEDINA CARABALLO"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# --- 2. Update the summary values in the header block. ---
# VALOR MORA total grows to reflect the new period being added.
$ws.Range("E11").Value = 113880

# Cant. Periodos goes from 1 to 2 (two periods of mora now on file).
$ws.Range("F13").Value = 2

# The existing 2507 row's "Salario Basico" total is also corrected.
$ws.Range("G16").Value = 1423500
